# Actualización automática del index.html y archivo Excel
# Elimina la fila 59 (caso 6236 - San Jose 1157), desplazando hacia
# arriba las filas siguientes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(59).Delete()
